# "lucia cambia xls de clases"
# Merge the "Aronow et al 2015" reference onto the same line as the
# "Barabas 2010" reference (separated by a comma) instead of its own line,
# in the reading-list cell on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("C2").Value = "Gerber  & Green 2012. FEDAI [Descarga] https://drive.google.com/drive/folders/14HDw0lx7v8cduNtj2XNvvZ5fm_lQ7Z6y?usp=sharing)`nBarabas 2010 [pdf](https://drive.google.com/file/d/15SqCaheQIA_Eg8Q6CxkkF5Gdt2dPdK1Y/view), Aronow et al 2015 [pdf](url)"

# After editing the cell, the active selection moved to C3.
$ws.Range("C3").Select()
